$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two remaining data rows with their new values
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "xuxsss"
$ws.Range("C2").Value = 1235

$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "zinhos"
$ws.Range("C3").Value = 232

# Remove rows 4-8, which are no longer part of the product list
$ws.Range("A4:D8").EntireRow.Delete()
